# Adds a new row (row 53) of sensor data to each of the four worksheets,
# matching the pattern of existing rows (time, lengths, checksum, decimal
# equivalents). Large integer-like values in column G are forced to text
# so Excel does not round them to floating point.

$wb = $excel.ActiveWorkbook

function Add-SensorRow {
    param($Worksheet, $RowIndex, $TimeStamp, $TotalLenHex, $IdHex, $ActualLenHex, $ChecksumHex, $TotalLenDec, $IdDec, $ActualLenDec, $ChecksumDec)

    $Worksheet.Range("A$RowIndex").Value = $TimeStamp
    $Worksheet.Range("B$RowIndex").Value = $TotalLenHex
    $Worksheet.Range("C$RowIndex").Value = $IdHex
    $Worksheet.Range("D$RowIndex").Value = $ActualLenHex
    $Worksheet.Range("E$RowIndex").Value = $ChecksumHex
    $Worksheet.Range("F$RowIndex").Value = $TotalLenDec

    # ID_DEC can exceed numeric precision (25-digit number), so store it
    # explicitly as text to preserve every digit.
    $gCell = $Worksheet.Range("G$RowIndex")
    $gCell.NumberFormat = "@"
    $gCell.Value = $IdDec

    $Worksheet.Range("H$RowIndex").Value = $ActualLenDec
    $Worksheet.Range("I$RowIndex").Value = $ChecksumDec
}

$wsFe35 = $wb.Worksheets.Item("ROW35-FE-LIFTER")
Add-SensorRow $wsFe35 53 "2025-03-06 12:42:06" "0x01,0x90 " "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c," "0x01,0x90," "0x d" 400 "568631262647113770877196" 400 13

$wsMid35 = $wb.Worksheets.Item("ROW35-MID-LIFTER")
Add-SensorRow $wsMid35 53 "2025-03-06 12:29:35" "0x01,0x90 " "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c," "0x01,0x90," "0x e" 400 "568631262647113770942732" 400 14

$wsFe02 = $wb.Worksheets.Item("ROW02-FE-LIFTER")
Add-SensorRow $wsFe02 53 "2025-03-06 12:51:45" "0x01,0x90 " "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x06,0x41,0x0c," "0x01,0x90," "0xff" 400 "568631262647113769959692" 400 255

$wsMid02 = $wb.Worksheets.Item("ROW02-MID-LIFTER")
Add-SensorRow $wsMid02 53 "2025-03-06 12:41:15" "0x01,0x90 " "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c," "0x01,0x90," "0x 3" 400 "568631262647113769959692" 400 3

Write-Host "Added row 53 to all four worksheets."
